$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9179281773574478
$ws.Range("J2").Value = 0.9179281773574478
$ws.Range("M2").Value = 15.25136533333333
$ws.Range("N2").Value = 45.754096
$ws.Range("O2").Value = 0.3045917506163436
$ws.Range("P2").Value = 0.3045917506163436
$ws.Range("Q2").Value = 9.406701523774222
$ws.Range("R2").Value = 84.66031371396799
$ws.Range("S2").Value = 0.2795933504813746
$ws.Range("T2").Value = 0.2795933504813746
$ws.Range("I3").Value = 0.9179281773574478
$ws.Range("J3").Value = 0.9179281773574478
$ws.Range("O3").Value = 0.2962340951184504
$ws.Range("P3").Value = 0.2962340951184504
$ws.Range("S3").Value = 0.271921623003212
$ws.Range("T3").Value = 0.271921623003212
$ws.Range("I4").Value = 0.9179281773574478
$ws.Range("J4").Value = 0.9179281773574478
$ws.Range("M4").Value = 11.41743966666667
$ws.Range("N4").Value = 34.252319
$ws.Range("O4").Value = 0.2280227284324326
$ws.Range("P4").Value = 0.2280227284324326
$ws.Range("Q4").Value = 7.042021796914111
$ws.Range("R4").Value = 63.378196172227
$ws.Range("S4").Value = 0.2093084875060551
$ws.Range("T4").Value = 0.2093084875060551
$ws.Range("I5").Value = 0.9179281773574478
$ws.Range("J5").Value = 0.9179281773574478
$ws.Range("M5").Value = 8.569808333333333
$ws.Range("N5").Value = 25.709425
$ws.Range("O5").Value = 0.1711514258327733
$ws.Range("P5").Value = 0.1711514258327733
$ws.Range("Q5").Value = 5.285666387613889
$ws.Range("R5").Value = 47.570997488525
$ws.Range("S5").Value = 0.157104716366806
$ws.Range("T5").Value = 0.157104716366806
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.055146
$ws.Range("H6").Value = 0.165438
$ws.Range("I6").Value = 0.08207182264255215
$ws.Range("J6").Value = 0.08207182264255215
$ws.Range("M6").Value = 15.25136533333333
$ws.Range("N6").Value = 45.754096
$ws.Range("O6").Value = 0.3045917506163436
$ws.Range("P6").Value = 0.3045917506163436
$ws.Range("Q6").Value = 0.8410517926719999
$ws.Range("R6").Value = 7.569466134048
$ws.Range("S6").Value = 0.02499840013496903
$ws.Range("T6").Value = 0.02499840013496903
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.055146
$ws.Range("H7").Value = 0.165438
$ws.Range("I7").Value = 0.08207182264255215
$ws.Range("J7").Value = 0.08207182264255215
$ws.Range("O7").Value = 0.2962340951184504
$ws.Range("P7").Value = 0.2962340951184504
$ws.Range("Q7").Value = 0.81797427621
$ws.Range("R7").Value = 7.36176848589
$ws.Range("S7").Value = 0.02431247211523839
$ws.Range("T7").Value = 0.02431247211523839
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.055146
$ws.Range("H8").Value = 0.165438
$ws.Range("I8").Value = 0.08207182264255215
$ws.Range("J8").Value = 0.08207182264255215
$ws.Range("M8").Value = 11.41743966666667
$ws.Range("N8").Value = 34.252319
$ws.Range("O8").Value = 0.2280227284324326
$ws.Range("P8").Value = 0.2280227284324326
$ws.Range("Q8").Value = 0.629626127858
$ws.Range("R8").Value = 5.666635150722
$ws.Range("S8").Value = 0.01871424092637744
$ws.Range("T8").Value = 0.01871424092637744
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.055146
$ws.Range("H9").Value = 0.165438
$ws.Range("I9").Value = 0.08207182264255215
$ws.Range("J9").Value = 0.08207182264255215
$ws.Range("M9").Value = 8.569808333333333
$ws.Range("N9").Value = 25.709425
$ws.Range("O9").Value = 0.1711514258327733
$ws.Range("P9").Value = 0.1711514258327733
$ws.Range("Q9").Value = 0.47259065035
$ws.Range("R9").Value = 4.25331585315
$ws.Range("S9").Value = 0.01404670946596729
$ws.Range("T9").Value = 0.01404670946596729
